$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @('Afsar Khan', 'REDX_VMS Portal', 'Mumbai Maharashtra India', 'test', 4.17, 0.06, 'Honeywell', '5MP (2560x1920)', 25, 'H265', 1, 30, 1),
    @('AK', 'REDX_VMS Portal', 'Kothrud, Pune, Maharashtra, India', 'Dallas', 4.17, 0.06, 'Honeywell', '5MP (2560x1920)', 25, 'H265', 1, 30, 1),
    @('AK', 'REDX_VMS Portal', 'Kothrud, Pune, Maharashtra, India', 'Dallas', 4.17, 0.06, 'Honeywell', '5MP (2560x1920)', 25, 'H265', 1, 30, 1),
    @('AK', 'REDX_VMS Portal', 'Kothrud, Pune, Maharashtra, India', 'Dallas', 4.17, 0.06, 'Honeywell', '5MP (2560x1920)', 25, 'H265', 1, 30, 100),
    @('AK', 'REDX_VMS Portal', 'Kothrud, Pune, Maharashtra, India', 'Dallas', 4.17, 0.06, 'Honeywell', '5MP (2560x1920)', 25, 'H265', 1, 30, 1),
    @('  ', '  ', '  ', '  ', 4.17, 0.06, 'Honeywell', '5MP (2560x1920)', 25, 'H265', 1, 30, 1)
)

$startRow = 342
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $vals = $rows[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $vals[$c]
    }
}
